# Adapt column header formatting to respective input file names.
# - Rename the "_old"/"_new" suffixed header labels to "_FV2310"/"_FV2404".
# - Turn the data range A1:U62 into a real Excel Table ("Table1").
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename header cells: "<name>_old" -> "<name>_FV2310" (A1:J1)
#    and "<name>_new" -> "<name>_FV2404" (L1:U1). K1 ("diff") is kept.
# ---------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$leftCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# ---------------------------------------------------------------------
# 2) Convert A1:U62 to a table without letting Excel snapshot the
#    existing header formatting into a new dxf (styles.xml must stay
#    untouched / dxfs count="0"). We stash the header formatting on a
#    scratch range outside the table, strip direct formatting from the
#    header row, add the table, then restore the formatting byte-for-
#    byte via copy/paste (which reuses the original style record
#    instead of synthesizing a header-row dxf override).
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A64:U64")

$headerRange.Copy($scratchRange)
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $null, 1)

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)
$scratchRange.Clear()

# ---------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
